# Append a new row (row 54) to each of the four worksheets, mirroring the
# last existing row (row 53) but stamped with the next day's timestamp.
# This expands each sheet's used range from A1:I53 to A1:I54.

$wb = $excel.ActiveWorkbook

$newTimestamp = 45840.46224537037

$rowsData = @{
    1 = @("0x01,0x90", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,", "0x01,0x68", "0x07", 400, "5.68631262647113e+23", 360, 7)
    2 = @("0x01,0x7c", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,", "0x01,0x60", "0x19", 380, "5.68432987514711e+23", 352, 25)
    3 = @("0x00,0x6e", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,", "0x00,0x69", "0x15", 110, "5.68631262647113e+23", 105, 15)
    4 = @("0x00,0x82", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,", "0x00,0x7E", "0x9", 130, "5.68631262647113e+23", 126, 9)
}

foreach ($sheetIndex in 1..4) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    $values = $rowsData[$sheetIndex]

    $sourceRow = 53
    $targetRow = 54

    # Column A: new timestamp, same date/time number format as the row above.
    $ws.Cells.Item($targetRow, 1).Value = $newTimestamp
    $ws.Cells.Item($targetRow, 1).NumberFormat = $ws.Cells.Item($sourceRow, 1).NumberFormat

    # Columns B-E: text-like hex byte strings.
    $ws.Cells.Item($targetRow, 2).Value = $values[0]
    $ws.Cells.Item($targetRow, 3).Value = $values[1]
    $ws.Cells.Item($targetRow, 4).Value = $values[2]
    $ws.Cells.Item($targetRow, 5).Value = $values[3]

    # Columns F-I: numeric decoded values.
    $ws.Cells.Item($targetRow, 6).Value = $values[4]
    $ws.Cells.Item($targetRow, 7).Value = [double]$values[5]
    $ws.Cells.Item($targetRow, 8).Value = $values[6]
    $ws.Cells.Item($targetRow, 9).Value = $values[7]
}
